$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 3).Value = 7318
}

for ($r = 36; $r -le 57; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}
